# ASSET ACCOUNTABILITY FORM - RETURN
# The "{#devices}{assignmentDate}" merge-field paragraph (first column of
# the repeating table row) needs an explicit 9pt (sz/szCs = 18 half-points)
# run size on each of its runs, matching the sizing already used elsewhere
# in the same table (e.g. the CONDITION column cells).

$d = $word.ActiveDocument

$target = $null
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "{#devices}{assignmentDate}*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # Setting both Size and SizeBi ensures the ASCII/Western run size (w:sz)
    # and the complex-script run size (w:szCs) are both stamped onto every
    # run in the paragraph.
    $r.Font.Size = 9
    $r.Font.SizeBi = 9
    Write-Host "Updated font size for merge-field paragraph."
} else {
    Write-Host "Target paragraph not found."
}
